$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.91
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 3.5
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("W2").Value = 12
$ws.Range("Z2").Value = 19
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 12
$ws.Range("AM2").Value = 23
$ws.Range("AY2").Value = 17
